# ADD results from server
# Update computed result values (row 2) on each year sheet.

$wb = $excel.ActiveWorkbook

$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("B2").Value = 348.111133040007
$ws2025.Range("E2").Value = 28878.40257356465
$ws2025.Range("I2").Value = 15751.17087451186
$ws2025.Range("L2").Value = 48991.24167597
$ws2025.Range("M2").Value = 11299.89730188
$ws2025.Range("N2").Value = 7337.013682751313
$ws2025.Range("O2").Value = 6991.647948766419

$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("A2").Value = 220.9181339077783
$ws2030.Range("B2").Value = 6050.092133217868
$ws2030.Range("E2").Value = 54981.0853406054
$ws2030.Range("I2").Value = 44136.25158196639
$ws2030.Range("L2").Value = 59075.72160612671
$ws2030.Range("M2").Value = 21956.38472896
$ws2030.Range("N2").Value = 10685.8390845415
$ws2030.Range("O2").Value = 9283.603099344651

$ws2035 = $wb.Worksheets.Item("2035")
$ws2035.Range("A2").Value = 2577.607803474762
$ws2035.Range("B2").Value = 7296.571711709315
$ws2035.Range("E2").Value = 65951.79560972707
$ws2035.Range("I2").Value = 59961.01146418095
$ws2035.Range("L2").Value = 59075.72160612671
$ws2035.Range("M2").Value = 28243.54178664774
$ws2035.Range("N2").Value = 15618.49026829105
$ws2035.Range("O2").Value = 15289.26274986199

$ws2040 = $wb.Worksheets.Item("2040")
$ws2040.Range("A2").Value = 2577.607803474762
$ws2040.Range("B2").Value = 7296.571711709315
$ws2040.Range("E2").Value = 65951.79560972707
$ws2040.Range("I2").Value = 59961.01146418095
$ws2040.Range("L2").Value = 59075.72160612671
$ws2040.Range("M2").Value = 28243.54178664774
$ws2040.Range("N2").Value = 15618.49026829105
$ws2040.Range("O2").Value = 15289.26274986199

$ws2045 = $wb.Worksheets.Item("2045")
$ws2045.Range("A2").Value = 2577.607803474762
$ws2045.Range("B2").Value = 7296.571711709315
$ws2045.Range("E2").Value = 65951.79560972707
$ws2045.Range("I2").Value = 59961.01146418095
$ws2045.Range("L2").Value = 59075.72160612671
$ws2045.Range("M2").Value = 28243.54178664774
$ws2045.Range("N2").Value = 15618.49026829105
$ws2045.Range("O2").Value = 15289.26274986199

$ws2050 = $wb.Worksheets.Item("2050")
$ws2050.Range("A2").Value = 2577.607803474762
$ws2050.Range("B2").Value = 7296.571711709315
$ws2050.Range("E2").Value = 65951.79560972707
$ws2050.Range("I2").Value = 59961.01146418095
$ws2050.Range("L2").Value = 59075.72160612671
$ws2050.Range("M2").Value = 28243.54178664774
$ws2050.Range("N2").Value = 15618.49026829105
$ws2050.Range("O2").Value = 15289.26274986199
